# Weekly update: insert a new price-record row for "Arándano (blue)" at
# Vega Monumental Concepción, pushing the existing history rows down by one.
#
# New dataset has a fresh week (2021-12-09) entered at the top of the
# detail rows (row 34); every previously existing row from 34..62 shifts
# down to 35..63 (dimension grows from A1:T62 to A1:T63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 34..62 down to 35..63, leaving a blank row 34
# (inherits formatting, e.g. the date number format on column D, from the
# row above - same as Excel's native Insert behaviour).
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with this week's record.
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44539
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101001
$ws.Range("J34").Value = "Arándano (blue)"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 290
$ws.Range("N34").Value = 3500
$ws.Range("O34").Value = 3600
$ws.Range("P34").Value = 3552
$ws.Range("Q34").Value = "$/bandeja 2 kilos"
$ws.Range("R34").Value = "Provincia de Linares"
$ws.Range("S34").Value = 1776
$ws.Range("T34").Value = 2
